$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.168.72"
$ws.Cells.Item(2, 5).Value = "  +1.09%  "
$ws.Cells.Item(3, 4).Value = "1.641.40"
$ws.Cells.Item(3, 5).Value = "  -0.07%  "
$ws.Cells.Item(4, 5).Value = "  +0.03%  "
$c = $ws.Cells.Item(5, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "217.19"
$c.Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  +0.20%  "
$ws.Cells.Item(6, 5).Value = "  +1.20%  "
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 5).Value = "  +0.77%  "
$c = $ws.Cells.Item(9, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0626"
$c.Style = $origStyle
$ws.Cells.Item(9, 5).Value = "  +0.71%  "
$c = $ws.Cells.Item(10, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "20.01"
$c.Style = $origStyle
$ws.Cells.Item(10, 5).Value = "  +0.84%  "
$c = $ws.Cells.Item(11, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0848"
$c.Style = $origStyle
$ws.Cells.Item(11, 5).Value = "  +0.22%  "
$ws.Cells.Item(12, 4).Value = "1.872.30"
$ws.Cells.Item(12, 5).Value = "  +0.01%  "
$ws.Cells.Item(13, 4).Value = "1.650.10"
$ws.Cells.Item(13, 5).Value = "  +0.51%  "
$c = $ws.Cells.Item(14, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.15"
$c.Style = $origStyle
$ws.Cells.Item(14, 5).Value = "  +0.47%  "
$c = $ws.Cells.Item(15, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.544"
$c.Style = $origStyle
$ws.Cells.Item(15, 5).Value = "  +2.75%  "
$c = $ws.Cells.Item(16, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "67.26"
$c.Style = $origStyle
$ws.Cells.Item(16, 5).Value = "  +1.23%  "
$ws.Cells.Item(17, 4).Value = "27.159.89"
$ws.Cells.Item(17, 5).Value = "  +1.03%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0739"
$ws.Cells.Item(18, 5).Value = "  +1.47%  "
$c = $ws.Cells.Item(19, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "218.46"
$c.Style = $origStyle
$ws.Cells.Item(19, 5).Value = "  -0.46%  "
$ws.Cells.Item(20, 5).Value = "  -0.01%  "
$ws.Cells.Item(21, 5).Value = "  +3.40%  "
$ws.Cells.Item(22, 5).Value = "  +6.53%  "
$ws.Cells.Item(23, 5).Value = "  +0.58%  "
$c = $ws.Cells.Item(24, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.21"
$c.Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  +0.26%  "
$c = $ws.Cells.Item(25, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "147.75"
$c.Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  +1.28%  "
$c = $ws.Cells.Item(26, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.53"
$c.Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  +1.79%  "
$ws.Cells.Item(27, 5).Value = "  -0.02%  "
$c = $ws.Cells.Item(28, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  -0.49%  "
$ws.Cells.Item(29, 5).Value = "  -0.25%  "
$ws.Cells.Item(30, 5).Value = "  +0.08%  "
$c = $ws.Cells.Item(31, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.19"
$c.Style = $origStyle
$ws.Cells.Item(31, 5).Value = "  +0.55%  "
$ws.Cells.Item(32, 5).Value = "  +0.68%  "
$c = $ws.Cells.Item(33, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.03"
$c.Style = $origStyle
$ws.Cells.Item(33, 5).Value = "  +0.97%  "
$ws.Cells.Item(34, 5).Value = "  +0.99%  "
$ws.Cells.Item(35, 4).Value = "1.273.00"
$ws.Cells.Item(35, 5).Value = "  +2.20%  "
$ws.Cells.Item(36, 5).Value = "  +0.76%  "
$c = $ws.Cells.Item(37, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0177"
$c.Style = $origStyle
$ws.Cells.Item(37, 5).Value = "  +1.89%  "
$c = $ws.Cells.Item(38, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.855"
$c.Style = $origStyle
$ws.Cells.Item(38, 5).Value = "  +2.58%  "
$ws.Cells.Item(39, 5).Value = "  +0.76%  "
$ws.Cells.Item(40, 5).Value = "  +0.01%  "
$ws.Cells.Item(41, 5).Value = "  -0.07%  "
$c = $ws.Cells.Item(42, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.Style = $origStyle
$ws.Cells.Item(42, 5).Value = "  +7.89%  "
$c = $ws.Cells.Item(43, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = $origStyle
$ws.Cells.Item(43, 5).Value = "  -1.18%  "
$ws.Cells.Item(44, 4).Value = "1.782.86"
$ws.Cells.Item(44, 5).Value = "  -0.03%  "
$c = $ws.Cells.Item(45, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "61.72"
$c.Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  +1.51%  "
$c = $ws.Cells.Item(46, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "91.83"
$c.Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  +0.32%  "
$ws.Cells.Item(47, 5).Value = "  +1.43%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Cells.Item(48, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0513"
$c.Style = $origStyle
$ws.Cells.Item(48, 5).Value = "  -0.23%  "
$ws.Cells.Item(49, 2).Value = "Algorand"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Cells.Item(49, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0974"
$c.Style = $origStyle
$ws.Cells.Item(49, 5).Value = "  -0.16%  "
$c = $ws.Cells.Item(50, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.63"
$c.Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  +0.52%  "
$ws.Cells.Item(51, 2).Value = "Mantle"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Cells.Item(51, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.405"
$c.Style = $origStyle
$ws.Cells.Item(51, 5).Value = "  +0.15%  "
